$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "x" marker in row 3 from column E (notes) to column D (output enable)
$ws.Range("E3").Value = $null
$ws.Range("D3").Value = "x"

# Update the active cell selection to C3
$ws.Range("C3").Select()
